$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The address list in column A (rows 4-27) had trailing newline characters
# baked into the cell text ("1488 4800 S\n", etc). Strip the stray
# trailing newline/carriage-return characters so each address is a clean
# single line.
for ($r = 4; $r -le 27; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val -ne $null) {
        $trimmed = $val.TrimEnd("`n", "`r")
        if ($trimmed -ne $val) {
            $cell.Value = $trimmed
        }
    }
}

# Shrink row 16's custom height slightly now that the wrapped text is
# shorter without the trailing newline.
$ws.Rows.Item(16).RowHeight = 36.55

# Move the active selection (used while checking the distance between the
# first and next addresses) to H20.
$ws.Range("H20").Select() | Out-Null
